$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1000
$ws.Range("I98").Value = 1000
$ws.Range("K98").Value = 1000
$ws.Range("M98").Value = 498

$ws.Range("H113").Value = 2851.8333
$ws.Range("I113").Value = 1805
$ws.Range("J113").Value = 3061.2
$ws.Range("K113").Value = 1805
$ws.Range("L113").Value = 3061.2
$ws.Range("M113").Value = 1449
$ws.Range("N113").Value = -9569.200000000001

$ws.Range("H122").Value = 1000
$ws.Range("I122").Value = 1000
$ws.Range("K122").Value = 3000
$ws.Range("M122").Value = -550

$ws.Range("H129").Value = 2060.625
$ws.Range("I129").Value = 297
$ws.Range("J129").Value = 2312.5715
$ws.Range("K129").Value = 891
$ws.Range("L129").Value = 6937.7145
$ws.Range("M129").Value = 4109
$ws.Range("N129").Value = -16937.7145

$ws.Range("H137").Value = 5053.25
$ws.Range("I137").Value = 1474
$ws.Range("J137").Value = 6246.3335
$ws.Range("K137").Value = 4422
$ws.Range("L137").Value = 18739.0005
$ws.Range("M137").Value = -1872
$ws.Range("N137").Value = -23839.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H101").Value = 28400.25
$ws.Range("J101").Value = 28400.25
$ws.Range("L101").Value = 28400.25
$ws.Range("N101").Value = -34890.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 43494.125
$ws.Range("I20").Value = 1697
$ws.Range("J20").Value = 127088.375
$ws.Range("K20").Value = 1697
$ws.Range("L20").Value = 127088.375
$ws.Range("M20").Value = -1450
$ws.Range("N20").Value = -127582.375

$ws.Range("H134").Value = 1584.7778
$ws.Range("I134").Value = 1511
$ws.Range("K134").Value = 4533
$ws.Range("M134").Value = -1998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1906.5714
$ws.Range("I99").Value = 1772.6666
$ws.Range("J99").Value = 2007
$ws.Range("K99").Value = 1772.6666
$ws.Range("L99").Value = 2007
$ws.Range("M99").Value = -274.6666
$ws.Range("N99").Value = -5003

$ws.Range("H126").Value = 1906.5714
$ws.Range("I126").Value = 1772.6666
$ws.Range("J126").Value = 2007
$ws.Range("K126").Value = 5317.9998
$ws.Range("L126").Value = 6021
$ws.Range("M126").Value = -2847.9998
$ws.Range("N126").Value = -10961

$ws.Range("H132").Value = 47627584
$ws.Range("I132").Value = 76934190
$ws.Range("J132").Value = 4344.75
$ws.Range("K132").Value = 230802570
$ws.Range("L132").Value = 13034.25
$ws.Range("M132").Value = -230800040
$ws.Range("N132").Value = -18094.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 170.5
$ws.Range("I50").Value = 141.28572
$ws.Range("J50").Value = 238.66667
$ws.Range("K50").Value = 423.85716
$ws.Range("L50").Value = 716.00001
$ws.Range("M50").Value = 57.14283999999998
$ws.Range("N50").Value = -1678.00001

$ws.Range("H53").Value = 170.5
$ws.Range("I53").Value = 141.28572
$ws.Range("J53").Value = 238.66667
$ws.Range("K53").Value = 423.85716
$ws.Range("L53").Value = 716.00001
$ws.Range("M53").Value = 57.14283999999998
$ws.Range("N53").Value = -1678.00001

$ws.Range("H87").Value = 5050
$ws.Range("I87").Value = 3575
$ws.Range("K87").Value = 10725
$ws.Range("M87").Value = -9477

$ws.Range("H90").Value = 5050
$ws.Range("I90").Value = 3575
$ws.Range("K90").Value = 32175
$ws.Range("M90").Value = -25935

$ws.Range("H131").Value = 672.30206
$ws.Range("I131").Value = 309.51282
$ws.Range("J131").Value = 920.5263
$ws.Range("K131").Value = 928.53846
$ws.Range("L131").Value = 2761.5789
$ws.Range("M131").Value = 4111.46154
$ws.Range("N131").Value = -12841.5789

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2366.1924
$ws.Range("I102").Value = 2131.5264
$ws.Range("J102").Value = 3003.1428
$ws.Range("K102").Value = 2131.5264
$ws.Range("L102").Value = 3003.1428
$ws.Range("M102").Value = -509.5264000000002
$ws.Range("N102").Value = -6247.1428

$ws.Range("H122").Value = 1592.5883
$ws.Range("I122").Value = 1589.3334
$ws.Range("J122").Value = 1596.25
$ws.Range("K122").Value = 4768.0002
$ws.Range("L122").Value = 4788.75
$ws.Range("M122").Value = -2318.0002
$ws.Range("N122").Value = -9688.75

$ws.Range("H123").Value = 23550.666
$ws.Range("J123").Value = 23550.666
$ws.Range("L123").Value = 23550.666
$ws.Range("N123").Value = -28450.666

$ws.Range("H126").Value = 1599.3334
$ws.Range("I126").Value = 1613.4286
$ws.Range("J126").Value = 1550
$ws.Range("K126").Value = 4840.2858
$ws.Range("L126").Value = 4650
$ws.Range("M126").Value = -2370.2858
$ws.Range("N126").Value = -9590

$ws.Range("H136").Value = 15703.75
$ws.Range("I136").Value = 10000
$ws.Range("J136").Value = 16518.572
$ws.Range("K136").Value = 30000
$ws.Range("L136").Value = 49555.716
$ws.Range("M136").Value = -27450
$ws.Range("N136").Value = -54655.716

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2165.818
$ws.Range("I7").Value = 1501
$ws.Range("J7").Value = 2545.7144
$ws.Range("K7").Value = 1501
$ws.Range("L7").Value = 2545.7144
$ws.Range("M7").Value = -1389
$ws.Range("N7").Value = -2769.7144

$ws.Range("H40").Value = 2415.2666
$ws.Range("I40").Value = 1530.4445
$ws.Range("K40").Value = 1530.4445
$ws.Range("M40").Value = -1394.4445

$ws.Range("H97").Value = 16448
$ws.Range("J97").Value = 16448
$ws.Range("L97").Value = 16448
$ws.Range("N97").Value = -18430

$ws.Range("H126").Value = 2165.818
$ws.Range("I126").Value = 1501
$ws.Range("J126").Value = 2545.7144
$ws.Range("K126").Value = 4503
$ws.Range("L126").Value = 7637.1432
$ws.Range("M126").Value = -2033
$ws.Range("N126").Value = -12577.1432

$ws.Range("H132").Value = 32374.5
$ws.Range("I132").Value = 58750
$ws.Range("J132").Value = 5999
$ws.Range("K132").Value = 176250
$ws.Range("L132").Value = 17997
$ws.Range("M132").Value = -173720
$ws.Range("N132").Value = -23057

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 2726.1428
$ws.Range("I62").Value = 2712
$ws.Range("K62").Value = 2712
$ws.Range("M62").Value = -2088

$ws.Range("H65").Value = 2726.1428
$ws.Range("I65").Value = 2712
$ws.Range("K65").Value = 13560
$ws.Range("M65").Value = -10440

$ws.Range("H106").Value = 27666.666
$ws.Range("J106").Value = 27666.666
$ws.Range("L106").Value = 27666.666
$ws.Range("N106").Value = -30190.666

$ws.Range("H122").Value = 13334852
$ws.Range("I122").Value = 22223556
$ws.Range("J122").Value = 1796.6666
$ws.Range("K122").Value = 66670668
$ws.Range("L122").Value = 5389.9998
$ws.Range("M122").Value = -66668218
$ws.Range("N122").Value = -10289.9998

$ws.Range("H126").Value = 1101
$ws.Range("I126").Value = 934.6667
$ws.Range("J126").Value = 1600
$ws.Range("K126").Value = 2804.0001
$ws.Range("L126").Value = 4800
$ws.Range("M126").Value = -334.0001000000002
$ws.Range("N126").Value = -9740

$ws.Range("H132").Value = 38472684
$ws.Range("I132").Value = 62516124
$ws.Range("J132").Value = 3180.6
$ws.Range("K132").Value = 187548372
$ws.Range("L132").Value = 9541.799999999999
$ws.Range("M132").Value = -187545842
$ws.Range("N132").Value = -14601.8
